$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3491903333333333
$ws.Range("H2").Value = 1.047571
$ws.Range("I2").Value = 0.008130334326258625
$ws.Range("J2").Value = 0.008130334326258625
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05968133333333333
$ws.Range("N2").Value = 0.179044
$ws.Range("O2").Value = 0.02602747651633847
$ws.Range("P2").Value = 0.02602747651633848
$ws.Range("Q2").Value = 0.02084014468044444
$ws.Range("R2").Value = 0.187561302124
$ws.Range("S2").Value = 0.0002116120857466769
$ws.Range("T2").Value = 0.000211612085746677

$ws.Range("G3").Value = 0.3491903333333333
$ws.Range("H3").Value = 1.047571
$ws.Range("I3").Value = 0.008130334326258625
$ws.Range("J3").Value = 0.008130334326258625
$ws.Range("O3").Value = 0.144012433133819
$ws.Range("P3").Value = 0.144012433133819
$ws.Range("Q3").Value = 0.1153104466507778
$ws.Range("R3").Value = 1.037794019857
$ws.Range("S3").Value = 0.001170869228515914
$ws.Range("T3").Value = 0.001170869228515914

$ws.Range("G4").Value = 0.3491903333333333
$ws.Range("H4").Value = 1.047571
$ws.Range("I4").Value = 0.008130334326258625
$ws.Range("J4").Value = 0.008130334326258625
$ws.Range("O4").Value = 0.8299600903498424
$ws.Range("P4").Value = 0.8299600903498425
$ws.Range("Q4").Value = 0.6645472660796666
$ws.Range("R4").Value = 5.980925394717
$ws.Range("S4").Value = 0.006747853011996034
$ws.Range("T4").Value = 0.006747853011996035

$ws.Range("I5").Value = 0.801301577139928
$ws.Range("J5").Value = 0.8013015771399279
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05968133333333333
$ws.Range("N5").Value = 0.179044
$ws.Range("O5").Value = 0.02602747651633847
$ws.Range("P5").Value = 0.02602747651633848
$ws.Range("Q5").Value = 2.053942695361334
$ws.Range("R5").Value = 18.485484258252
$ws.Range("S5").Value = 0.02085585798151446
$ws.Range("T5").Value = 0.02085585798151446

$ws.Range("I6").Value = 0.801301577139928
$ws.Range("J6").Value = 0.8013015771399279
$ws.Range("O6").Value = 0.144012433133819
$ws.Range("P6").Value = 0.144012433133819
$ws.Range("S6").Value = 0.1153973897978876
$ws.Range("T6").Value = 0.1153973897978876

$ws.Range("I7").Value = 0.801301577139928
$ws.Range("J7").Value = 0.8013015771399279
$ws.Range("O7").Value = 0.8299600903498424
$ws.Range("P7").Value = 0.8299600903498425
$ws.Range("S7").Value = 0.6650483293605258
$ws.Range("T7").Value = 0.6650483293605258

$ws.Range("I8").Value = 0.1905680885338134
$ws.Range("J8").Value = 0.1905680885338134
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05968133333333333
$ws.Range("N8").Value = 0.179044
$ws.Range("O8").Value = 0.02602747651633847
$ws.Range("P8").Value = 0.02602747651633848
$ws.Range("Q8").Value = 0.4884751816039999
$ws.Range("R8").Value = 4.396276634436
$ws.Range("S8").Value = 0.004960006449077339
$ws.Range("T8").Value = 0.00496000644907734

$ws.Range("I9").Value = 0.1905680885338134
$ws.Range("J9").Value = 0.1905680885338134
$ws.Range("O9").Value = 0.144012433133819
$ws.Range("P9").Value = 0.144012433133819
$ws.Range("S9").Value = 0.02744417410741551
$ws.Range("T9").Value = 0.02744417410741551

$ws.Range("I10").Value = 0.1905680885338134
$ws.Range("J10").Value = 0.1905680885338134
$ws.Range("O10").Value = 0.8299600903498424
$ws.Range("P10").Value = 0.8299600903498425
$ws.Range("S10").Value = 0.1581639079773206
$ws.Range("T10").Value = 0.1581639079773206
